$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 10 first (table shrinks from 9 data rows to 8 data rows)
$ws.Rows.Item(10).Delete()

# Overwrite data rows 2-9 with the refreshed TPM-derived values
# Row 2: FAPs -> ECs
$ws.Range("A2").Value2 = "FAPs"
$ws.Range("B2").Value2 = "Wnt5a"
$ws.Range("C2").Value2 = "Ror1"
$ws.Range("D2").Value2 = "ECs"
$ws.Range("E2").Value2 = 3
$ws.Range("F2").Value2 = 1
$ws.Range("G2").Value2 = 8.775006
$ws.Range("H2").Value2 = 26.325018
$ws.Range("I2").Value2 = 0.9920592728348052
$ws.Range("J2").Value2 = 0.9920592728348053
$ws.Range("K2").Value2 = 3
$ws.Range("L2").Value2 = 1
$ws.Range("M2").Value2 = 0.9943956666666667
$ws.Range("N2").Value2 = 2.983187
$ws.Range("O2").Value2 = 0.03092298537432404
$ws.Range("P2").Value2 = 0.03092298537432405
$ws.Range("Q2").Value2 = 8.725827941374
$ws.Range("R2").Value2 = 78.532451472366
$ws.Range("S2").Value2 = 0.03067743438433323
$ws.Range("T2").Value2 = 0.03067743438433324

# Row 3: FAPs -> FAPs
$ws.Range("A3").Value2 = "FAPs"
$ws.Range("B3").Value2 = "Wnt5a"
$ws.Range("C3").Value2 = "Ror1"
$ws.Range("D3").Value2 = "FAPs"
$ws.Range("E3").Value2 = 3
$ws.Range("F3").Value2 = 1
$ws.Range("G3").Value2 = 8.775006
$ws.Range("H3").Value2 = 26.325018
$ws.Range("I3").Value2 = 0.9920592728348052
$ws.Range("J3").Value2 = 0.9920592728348053
$ws.Range("K3").Value2 = 3
$ws.Range("L3").Value2 = 1
$ws.Range("M3").Value2 = 13.33283666666667
$ws.Range("N3").Value2 = 39.99851
$ws.Range("O3").Value2 = 0.4146147525196222
$ws.Range("P3").Value2 = 0.4146147525196222
$ws.Range("Q3").Value2 = 116.99572174702
$ws.Range("R3").Value2 = 1052.96149572318
$ws.Range("S3").Value2 = 0.4113224098911991
$ws.Range("T3").Value2 = 0.4113224098911992

# Row 4: FAPs -> MuSCs
$ws.Range("A4").Value2 = "FAPs"
$ws.Range("B4").Value2 = "Wnt5a"
$ws.Range("C4").Value2 = "Ror1"
$ws.Range("D4").Value2 = "MuSCs"
$ws.Range("E4").Value2 = 3
$ws.Range("F4").Value2 = 1
$ws.Range("G4").Value2 = 8.775006
$ws.Range("H4").Value2 = 26.325018
$ws.Range("I4").Value2 = 0.9920592728348052
$ws.Range("J4").Value2 = 0.9920592728348053
$ws.Range("K4").Value2 = 3
$ws.Range("L4").Value2 = 1
$ws.Range("M4").Value2 = 17.81414066666667
$ws.Range("N4").Value2 = 53.442422
$ws.Range("O4").Value2 = 0.553971049711082
$ws.Range("P4").Value2 = 0.553971049711082
$ws.Range("Q4").Value2 = 156.319191234844
$ws.Range("R4").Value2 = 1406.872721113596
$ws.Range("S4").Value2 = 0.5495721167479097
$ws.Range("T4").Value2 = 0.5495721167479098

# Row 5: FAPs -> Resolving-Mac
$ws.Range("A5").Value2 = "FAPs"
$ws.Range("B5").Value2 = "Wnt5a"
$ws.Range("C5").Value2 = "Ror1"
$ws.Range("D5").Value2 = "Resolving-Mac"
$ws.Range("E5").Value2 = 3
$ws.Range("F5").Value2 = 1
$ws.Range("G5").Value2 = 8.775006
$ws.Range("H5").Value2 = 26.325018
$ws.Range("I5").Value2 = 0.9920592728348052
$ws.Range("J5").Value2 = 0.9920592728348053
$ws.Range("K5").Value2 = 2
$ws.Range("L5").Value2 = 0.6666666666666666
$ws.Range("M5").Value2 = 0.015796
$ws.Range("N5").Value2 = 0.047388
$ws.Range("O5").Value2 = 0.0004912123949717091
$ws.Range("P5").Value2 = 0.0004912123949717091
$ws.Range("Q5").Value2 = 0.138609994776
$ws.Range("R5").Value2 = 1.247489952984
$ws.Range("S5").Value2 = 0.0004873118113630769
$ws.Range("T5").Value2 = 0.0004873118113630769

# Row 6: MuSCs -> ECs
$ws.Range("A6").Value2 = "MuSCs"
$ws.Range("B6").Value2 = "Wnt5a"
$ws.Range("C6").Value2 = "Ror1"
$ws.Range("D6").Value2 = "ECs"
$ws.Range("E6").Value2 = 1
$ws.Range("F6").Value2 = 0.3333333333333333
$ws.Range("G6").Value2 = 0.07023766666666667
$ws.Range("H6").Value2 = 0.210713
$ws.Range("I6").Value2 = 0.007940727165194733
$ws.Range("J6").Value2 = 0.007940727165194734
$ws.Range("K6").Value2 = 3
$ws.Range("L6").Value2 = 1
$ws.Range("M6").Value2 = 0.9943956666666667
$ws.Range("N6").Value2 = 2.983187
$ws.Range("O6").Value2 = 0.03092298537432404
$ws.Range("P6").Value2 = 0.03092298537432405
$ws.Range("Q6").Value2 = 0.06984403137011111
$ws.Range("R6").Value2 = 0.6285962823310001
$ws.Range("S6").Value2 = 0.0002455509899908143
$ws.Range("T6").Value2 = 0.0002455509899908144

# Row 7: MuSCs -> FAPs
$ws.Range("A7").Value2 = "MuSCs"
$ws.Range("B7").Value2 = "Wnt5a"
$ws.Range("C7").Value2 = "Ror1"
$ws.Range("D7").Value2 = "FAPs"
$ws.Range("E7").Value2 = 1
$ws.Range("F7").Value2 = 0.3333333333333333
$ws.Range("G7").Value2 = 0.07023766666666667
$ws.Range("H7").Value2 = 0.210713
$ws.Range("I7").Value2 = 0.007940727165194733
$ws.Range("J7").Value2 = 0.007940727165194734
$ws.Range("K7").Value2 = 3
$ws.Range("L7").Value2 = 1
$ws.Range("M7").Value2 = 13.33283666666667
$ws.Range("N7").Value2 = 39.99851
$ws.Range("O7").Value2 = 0.4146147525196222
$ws.Range("P7").Value2 = 0.4146147525196222
$ws.Range("Q7").Value2 = 0.9364673375144446
$ws.Range("R7").Value2 = 8.428206037630002
$ws.Range("S7").Value2 = 0.003292342628423055
$ws.Range("T7").Value2 = 0.003292342628423056

# Row 8: MuSCs -> MuSCs
$ws.Range("A8").Value2 = "MuSCs"
$ws.Range("B8").Value2 = "Wnt5a"
$ws.Range("C8").Value2 = "Ror1"
$ws.Range("D8").Value2 = "MuSCs"
$ws.Range("E8").Value2 = 1
$ws.Range("F8").Value2 = 0.3333333333333333
$ws.Range("G8").Value2 = 0.07023766666666667
$ws.Range("H8").Value2 = 0.210713
$ws.Range("I8").Value2 = 0.007940727165194733
$ws.Range("J8").Value2 = 0.007940727165194734
$ws.Range("K8").Value2 = 3
$ws.Range("L8").Value2 = 1
$ws.Range("M8").Value2 = 17.81414066666667
$ws.Range("N8").Value2 = 53.442422
$ws.Range("O8").Value2 = 0.553971049711082
$ws.Range("P8").Value2 = 0.553971049711082
$ws.Range("Q8").Value2 = 1.251223674098445
$ws.Range("R8").Value2 = 11.261013066886
$ws.Range("S8").Value2 = 0.004398932963172231
$ws.Range("T8").Value2 = 0.004398932963172232

# Row 9: MuSCs -> Resolving-Mac
$ws.Range("A9").Value2 = "MuSCs"
$ws.Range("B9").Value2 = "Wnt5a"
$ws.Range("C9").Value2 = "Ror1"
$ws.Range("D9").Value2 = "Resolving-Mac"
$ws.Range("E9").Value2 = 1
$ws.Range("F9").Value2 = 0.3333333333333333
$ws.Range("G9").Value2 = 0.07023766666666667
$ws.Range("H9").Value2 = 0.210713
$ws.Range("I9").Value2 = 0.007940727165194733
$ws.Range("J9").Value2 = 0.007940727165194734
$ws.Range("K9").Value2 = 2
$ws.Range("L9").Value2 = 0.6666666666666666
$ws.Range("M9").Value2 = 0.015796
$ws.Range("N9").Value2 = 0.047388
$ws.Range("O9").Value2 = 0.0004912123949717091
$ws.Range("P9").Value2 = 0.0004912123949717091
$ws.Range("Q9").Value2 = 0.001109474182666667
$ws.Range("R9").Value2 = 0.009985267644
$ws.Range("S9").Value2 = (3.900583608632215 / 1000000.0)
$ws.Range("T9").Value2 = (3.900583608632215 / 1000000.0)

Write-Output "Edit complete"